$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.598999999999999
$ws.Range("E5").Value = 12.94
$ws.Range("E9").Value = 13.152
$ws.Range("E11").Value = 13.246
$ws.Range("B21").Value = 6.311999999999999
$ws.Range("E21").Value = 12.694
$ws.Range("B23").Value = 6.842000000000001
$ws.Range("B25").Value = 6.439
